# MHV-17222: bump the CodeSystem "Version" and "Date" metadata values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 0.2.9-beta -> 0.2.10-beta  (row 3, column B)
$ws.Range("B3").Value = "0.2.10-beta"

# Date: 2023-02-16T09:21:54-06:00 -> 2023-12-06T12:46:33-06:00  (row 8, column B)
$ws.Range("B8").Value = "2023-12-06T12:46:33-06:00"
